$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14 model name (first new shared string)
$ws.Range("A14").Value = "Soil Ag~Chem + plants,  chemistry~microbes, microbes~plants"

# New header cells in row 1 (CFI, SRMR)
$ws.Range("O1").Value = "CFI"
$ws.Range("P1").Value = "SRMR"

# New values in row 8 (CFI, SRMR for that model)
$ws.Range("O8").Value = 0.951
$ws.Range("P8").Value = 0.088

# Remaining new row 14 data
$ws.Range("E14").Value = 0.075
$ws.Range("O14").Value = 0.95
$ws.Range("P14").Value = 0.087
$ws.Range("Q14").Value = "(NOT SIG)"

# Update selection to match new active cell
$ws.Range("Q14").Select()
